# MarsUndiscovered/Spreadsheets/Breeds.xlsx
#
# Add a "line attack" command: two new breed-stat columns, LineAttackMin
# and LineAttackMax, holding the min/max damage for an attack that hits
# every object along a path (as opposed to the existing single-target
# lightning/basic attacks). Only the header row is populated for now -
# no breed currently uses the new command, so rows 2-4 stay blank in K/L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Breeds")

$ws.Range("K1").Value = "LineAttackMin"
$ws.Range("L1").Value = "LineAttackMax"

# Leave the sheet scrolled/selected the way it was after adding the
# columns: column D at the left edge of the view, with the cell just
# past the new data (M1) selected.
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M1").Select()
